$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (shared strings) with shortened names ---
$ws.Range("C1").Value = "Num_method"
$ws.Range("D1").Value = "Mesh_dim"
$ws.Range("F1").Value = "Bound_cond"
$ws.Range("J1").Value = "Comput_time"

# --- Rewrite data rows 2-14 with the synthesized/updated values ---
# Row 2 (case id 2)
$ws.Range("A2").Value = 2
$ws.Range("C2").Value = "FE"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Regular_RightTriangles"
$ws.Range("F2").Value = "Dirichlet"
$ws.Range("G2").Value = 2.0039
$ws.Range("H2").Value = "Triangles"
$ws.Range("I2").Value = "Green"
$ws.Range("J2").Value = 108.1455399990082

# Row 3 (case id 3)
$ws.Range("A3").Value = 3
$ws.Range("C3").Value = "FE"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Unstructured_triangles"
$ws.Range("F3").Value = "Dirichlet"
$ws.Range("G3").Value = 2.0156
$ws.Range("H3").Value = "Triangles"
$ws.Range("I3").Value = "Green"
$ws.Range("J3").Value = 6.762243032455444

# Row 4 (case id 6)
$ws.Range("A4").Value = 6
$ws.Range("C4").Value = "FE"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "Regular_Tetrahedra"
$ws.Range("F4").Value = "Dirichlet"
$ws.Range("G4").Value = 1.3403
$ws.Range("H4").Value = "Tetrahedron"
$ws.Range("I4").Value = "Green"
$ws.Range("J4").Value = 210.4683861732483

# Row 5 (case id 7)
$ws.Range("A5").Value = 7
$ws.Range("C5").Value = "FE"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "Unstructured_Tetrahedra"
$ws.Range("F5").Value = "Dirichlet"
$ws.Range("G5").Value = 0.6691
$ws.Range("H5").Value = "Tetrahedron"
$ws.Range("I5").Value = "Green"
$ws.Range("J5").Value = 11.9149010181427

# Row 6 (case id 1)
$ws.Range("A6").Value = 1
$ws.Range("C6").Value = "FV"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "RegularSquares"
$ws.Range("F6").Value = "Dirichlet"
$ws.Range("G6").Value = 2.0039
$ws.Range("H6").Value = "Squares"
$ws.Range("I6").Value = "Green"
$ws.Range("J6").Value = 9.832487106323242

# Row 7 (case id 11)
$ws.Range("A7").Value = 11
$ws.Range("C7").Value = "FV"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "RegularSquares"
$ws.Range("F7").Value = "Neumann"
$ws.Range("G7").Value = 2.0039
$ws.Range("H7").Value = "Squares"
$ws.Range("I7").Value = "Green"
$ws.Range("J7").Value = 9.859630107879639

# Row 8 (case id 4)
$ws.Range("A8").Value = 4
$ws.Range("C8").Value = "FV"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Regular_RightTriangles"
$ws.Range("F8").Value = "Dirichlet"
$ws.Range("G8").Value = 0.0212
$ws.Range("H8").Value = "Triangles"
$ws.Range("I8").Value = "Green"
$ws.Range("J8").Value = 15.65501999855042

# Row 9 (case id 0)
$ws.Range("A9").Value = 0
$ws.Range("C9").Value = "FV"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Regular_RightTriangles"
$ws.Range("F9").Value = "Neumann"
$ws.Range("G9").Value = -0.0056
$ws.Range("H9").Value = "Triangles"
$ws.Range("I9").Value = "Orange(order 0)"
$ws.Range("J9").Value = 15.78992199897766

# Row 10 (case id 8)
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = "FV"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Unstructured_triangles"
$ws.Range("F10").Value = "Dirichlet"
$ws.Range("G10").Value = 0.6138
$ws.Range("H10").Value = "Triangles"
$ws.Range("I10").Value = "Green"
$ws.Range("J10").Value = 2.600184917449951

# Row 11 (case id 10)
$ws.Range("A11").Value = 10
$ws.Range("C11").Value = "FV"
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = "Regular_Cubes"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = 1.3403
$ws.Range("H11").Value = "Cubes"
$ws.Range("I11").Value = "Green"
$ws.Range("J11").Value = 5.900697946548462

# Row 12 (case id 9)
$ws.Range("A12").Value = 9
$ws.Range("C12").Value = "FV"
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "Regular_Tetrahedra"
$ws.Range("F12").Value = "Dirichlet"
$ws.Range("G12").Value = 0.0065
$ws.Range("H12").Value = "Tetrahedron"
$ws.Range("I12").Value = "Green"
$ws.Range("J12").Value = 62.56098890304565

# Row 13 (case id 12)
$ws.Range("A13").Value = 12
$ws.Range("C13").Value = "FV"
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "Unstructured_Tetrahedra"
$ws.Range("F13").Value = "Dirichlet"
$ws.Range("G13").Value = 0.5359
$ws.Range("H13").Value = "Tetrahedron"
$ws.Range("I13").Value = "Green"
$ws.Range("J13").Value = 3.782500028610229

# Row 14 (case id 5)
$ws.Range("A14").Value = 5
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "Structured_triangles"
$ws.Range("F14").Value = "Dirichlet"
$ws.Range("G14").Value = 0.8952
$ws.Range("H14").Value = "Triangles"
$ws.Range("I14").Value = "Green"
$ws.Range("J14").Value = 4.790747165679932
